# Apply the scenarios.xlsx update:
#  - SolverSettings: enable include_growth_limit (N -> Y) for all scenarios
#  - PowerPlants: drop EC_COAL / EC_DSL_CC / EC_OIL and EI_BATT / EI_SOLPV rows,
#    add a new ED_NG_CC row, and refresh which scenarios have each asset active
#  - Fuels: BIO no longer active for scenario A

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: SolverSettings
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SolverSettings")
$ws1.Range("B6:E6").Value = "Y"

# ---------------------------------------------------------------------------
# Sheet 2: PowerPlants
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("PowerPlants")

# Remove the retired capacity-addition options (EC_COAL, EC_DSL_CC, EC_OIL)
$ws2.Range("A15:A17").EntireRow.Delete()

# Remove the inactive storage/solar "EI_" rows (now at rows 24:25)
$ws2.Range("A24:A25").EntireRow.Delete()

# Insert a new row for the ED_NG_CC retirement option, ahead of ED_NG_OC
$ws2.Rows(21).Insert()
$ws2.Range("A21").Value = "ED_NG_CC"

# Refresh the Y markers for the EC_*/ED_* rows (13-24)
$ws2.Range("B13").Value = "Y"
$ws2.Range("C13").Value = "Y"
$ws2.Range("D13").ClearContents()
$ws2.Range("E13").ClearContents()

$ws2.Range("B14").ClearContents()
$ws2.Range("C14").Value = "Y"
$ws2.Range("D14").ClearContents()
$ws2.Range("E14").ClearContents()

$ws2.Range("B15").Value = "Y"
$ws2.Range("C15").Value = "Y"
$ws2.Range("D15").ClearContents()
$ws2.Range("E15").ClearContents()

$ws2.Range("B16").Value = "Y"
$ws2.Range("C16").Value = "Y"
$ws2.Range("D16").ClearContents()
$ws2.Range("E16").ClearContents()

$ws2.Range("B17").ClearContents()
$ws2.Range("C17").Value = "Y"
$ws2.Range("D17").ClearContents()
$ws2.Range("E17").ClearContents()

$ws2.Range("B18").ClearContents()
$ws2.Range("C18").Value = "Y"
$ws2.Range("D18").ClearContents()
$ws2.Range("E18").ClearContents()

$ws2.Range("B19").ClearContents()
$ws2.Range("C19").ClearContents()
$ws2.Range("D19").Value = "Y"
$ws2.Range("E19").Value = "Y"

$ws2.Range("B20").ClearContents()
$ws2.Range("C20").ClearContents()
$ws2.Range("D20").ClearContents()
$ws2.Range("E20").Value = "Y"

$ws2.Range("B21").ClearContents()
$ws2.Range("C21").ClearContents()
$ws2.Range("D21").Value = "Y"
$ws2.Range("E21").Value = "Y"

$ws2.Range("B22").ClearContents()
$ws2.Range("C22").ClearContents()
$ws2.Range("D22").Value = "Y"
$ws2.Range("E22").Value = "Y"

$ws2.Range("B23").ClearContents()
$ws2.Range("C23").ClearContents()
$ws2.Range("D23").ClearContents()
$ws2.Range("E23").Value = "Y"

$ws2.Range("B24").ClearContents()
$ws2.Range("C24").ClearContents()
$ws2.Range("D24").ClearContents()
$ws2.Range("E24").Value = "Y"

# ---------------------------------------------------------------------------
# Sheet 3: Fuels
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Fuels")
$ws3.Range("B2").ClearContents()
$ws3.Range("B7").Select()

# ---------------------------------------------------------------------------
# Restore selections to match the final saved view for each sheet
# ---------------------------------------------------------------------------
$ws2.Range("D22").Select()
$ws1.Activate()
$ws1.Range("E8").Select()
